$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (division problems row #1): 88÷9=, 29÷9=, 96÷2=, 82÷2=, 91÷7=
$row = $t.Rows.Item(1)
$row.Cells.Item(1).Range.Text = "36÷7="
$row.Cells.Item(2).Range.Text = "71÷4="
$row.Cells.Item(3).Range.Text = "92÷9="
$row.Cells.Item(4).Range.Text = "79÷5="
$row.Cells.Item(5).Range.Text = "24÷5="

# Row 5 (division problems row #2): 21÷6=, 83÷4=, 76÷7=, 99÷9=, 53÷5=
$row = $t.Rows.Item(5)
$row.Cells.Item(1).Range.Text = "69÷4="
$row.Cells.Item(2).Range.Text = "22÷4="
$row.Cells.Item(3).Range.Text = "30÷9="
$row.Cells.Item(4).Range.Text = "44÷4="
$row.Cells.Item(5).Range.Text = "10÷2="

# Row 9 (division problems row #3): 18÷7=, 25÷4=, 25÷6=, 37÷4=, 55÷3=
$row = $t.Rows.Item(9)
$row.Cells.Item(1).Range.Text = "26÷5="
$row.Cells.Item(2).Range.Text = "13÷4="
$row.Cells.Item(3).Range.Text = "79÷8="
$row.Cells.Item(4).Range.Text = "40÷8="
$row.Cells.Item(5).Range.Text = "28÷4="

# Row 13 (division problems row #4): 17÷6=, 26÷9=, 40÷3=, 78÷2=, 78÷3=
$row = $t.Rows.Item(13)
$row.Cells.Item(1).Range.Text = "21÷5="
$row.Cells.Item(2).Range.Text = "25÷2="
$row.Cells.Item(3).Range.Text = "73÷3="
$row.Cells.Item(4).Range.Text = "37÷6="
$row.Cells.Item(5).Range.Text = "26÷9="

# Row 17 (division problems row #5): 37÷5=, 98÷7=, 59÷2=, 78÷9=, 20÷6=
$row = $t.Rows.Item(17)
$row.Cells.Item(1).Range.Text = "27÷9="
$row.Cells.Item(2).Range.Text = "16÷9="
$row.Cells.Item(3).Range.Text = "43÷3="
$row.Cells.Item(4).Range.Text = "79÷4="
$row.Cells.Item(5).Range.Text = "88÷5="
